$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: 08-03-2017 / 15:28:00 / 25.5 / 17.75
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "08-03-2017"
$ws.Range("B5").Value = "15:28:00"
$ws.Range("C5").Value = 25.5
$ws.Range("D5").Value = 17.75

# Row 6: 12-03-2017 / 16:27:59 / 55.0 / 15.5
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "12-03-2017"
$ws.Range("B6").Value = "16:27:59"
$ws.Range("C6").Value = 55.0
$ws.Range("D6").Value = 15.5

# Row 7: 12-03-2017 / 18:47:34 / 56.0 / 16.0
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "12-03-2017"
$ws.Range("B7").Value = "18:47:34"
$ws.Range("C7").Value = 56.0
$ws.Range("D7").Value = 16.0

# Row 8: 12-03-2017 / 18:47:34 / 56.0 / 16.0
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "12-03-2017"
$ws.Range("B8").Value = "18:47:34"
$ws.Range("C8").Value = 56.0
$ws.Range("D8").Value = 16.0

# The dates in column A look like valid dates to Excel's auto-detection, so a
# text NumberFormat was applied above to force them to stay as text/strings
# (matching the original file's convention of storing dates as plain shared
# strings). Clear that temporary formatting again now that the values are
# safely stored as text, so no visible/number formatting is left behind.
$ws.Range("A5:A8").ClearFormats()
